# Regenerate the "K" column (column G) of save_data with newly computed
# strike/K values (replacing the old "Strike#" based values), as part of
# regenerating std/mean and writing s_vals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values for rows 2 through 40 (column G), in row order.
$newK = @(3, 3, 1, 5, 6, 7, 5, 8, 3, 9, 4, 2, 2, 8, 1, 7, 3, 3, 3, 5, 1, 7, 4, 6, 5, 7, 4, 8, 8, 5, 3, 3, 5, 7, 1, 7, 3, 2, 5)

for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
